# Generate Report for Handoff
#
# Updates the localization-status report with refreshed handoff timestamps
# and marks the affected rows' Priority column with the "ht" (handoff type)
# value produced by the latest handoff run.

$wb = $excel.ActiveWorkbook

$rows = 7,9,10,11,12,14

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ----------
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-24 20:22:17"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (column H) + Priority (E) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-24 20:22:10"
}

# --- de-de sheet: "Latest Handoff Datetime" (column H) + Priority (E) ----
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-24 20:22:17"
}
